$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 408, shifting rows 408:432 down to 409:433
$ws.Rows("408:408").Insert()

# Populate the newly inserted row 408 with the new weekly data entry
$ws.Cells.Item(408, 1).Value = 5
$ws.Cells.Item(408, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(408, 3).Value = "Maule"
$ws.Cells.Item(408, 4).Value = 45021
$ws.Cells.Item(408, 5).Value = 7
$ws.Cells.Item(408, 6).Value = 100112006
$ws.Cells.Item(408, 7).Value = "Repollo"
$ws.Cells.Item(408, 8).Value = "Crespo record"
$ws.Cells.Item(408, 9).Value = "Primera"
$ws.Cells.Item(408, 10).Value = 3000
$ws.Cells.Item(408, 11).Value = 1000
$ws.Cells.Item(408, 12).Value = 1000
$ws.Cells.Item(408, 13).Value = 1000
$ws.Cells.Item(408, 14).Value = "$/unidad"
$ws.Cells.Item(408, 15).Value = "Región del Maule"
$ws.Cells.Item(408, 16).Value = 1000
$ws.Cells.Item(408, 17).Value = 1
$ws.Cells.Item(408, 18).Value = "Hortaliza"
